$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update 想去人数 (F) values for several rows
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3376
$ws1.Range("F4").Value = 132
$ws1.Range("F5").Value = 6954
$ws1.Range("F6").Value = 2372
$ws1.Range("F7").Value = 36
$ws1.Range("F15").Value = 42

# Sheet "全部类型" (sheet4): same updates, rows shifted by one
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3376
$ws4.Range("F5").Value = 132
$ws4.Range("F6").Value = 6954
$ws4.Range("F7").Value = 2372
$ws4.Range("F8").Value = 36
$ws4.Range("F16").Value = 42
